$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 17.88092345040183
$ws.Range("D2").Value = 7.627201445290735
$ws.Range("E2").Value = 12.40420050516708
$ws.Range("F2").Value = 39.49635972393731
$ws.Range("G2").Value = 3.706666989924458
$ws.Range("K2").Value = 14.52246215023877
$ws.Range("L2").Value = 10.02768412981049
$ws.Range("N2").Value = 22.96549459106871

$ws.Range("B3").Value = 17.78584503217421
$ws.Range("D3").Value = 7.643603733329482
$ws.Range("E3").Value = 12.39044521242943
$ws.Range("F3").Value = 39.24893057530105
$ws.Range("G3").Value = 3.710291553195901
$ws.Range("K3").Value = 14.25318853072429
$ws.Range("L3").Value = 9.993460813934059
$ws.Range("N3").Value = 23.00150606009767

$ws.Range("B4").Value = 17.73153056602055
$ws.Range("D4").Value = 7.654150267535102
$ws.Range("E4").Value = 12.38192144256185
$ws.Range("F4").Value = 39.10596306349509
$ws.Range("G4").Value = 3.712632551796228
$ws.Range("K4").Value = 14.08999063151655
$ws.Range("L4").Value = 9.974723714928633
$ws.Range("N4").Value = 23.02547100231828

$ws.Range("B5").Value = 17.71043694566734
$ws.Range("D5").Value = 7.658567869174207
$ws.Range("E5").Value = 12.37842776358073
$ws.Range("F5").Value = 39.0499962951766
$ws.Range("G5").Value = 3.713615680382824
$ws.Range("K5").Value = 14.02412243908537
$ws.Range("L5").Value = 9.967665645246745
$ws.Range("N5").Value = 23.03570274434314

$ws.Range("B6").Value = 17.70699767322948
$ws.Range("D6").Value = 7.659308651624899
$ws.Range("E6").Value = 12.37784640586532
$ws.Range("F6").Value = 39.04084273139661
$ws.Range("G6").Value = 3.713780692076719
$ws.Range("K6").Value = 14.01322639002018
$ws.Range("L6").Value = 9.966528670327815
$ws.Range("N6").Value = 23.03742984359146

$ws.Range("B7").Value = 17.73124185668905
$ws.Range("D7").Value = 7.654209359473994
$ws.Range("E7").Value = 12.38187440795777
$ws.Range("F7").Value = 39.10519893704044
$ws.Range("G7").Value = 3.712645692431436
$ws.Range("K7").Value = 14.08909960334326
$ws.Range("L7").Value = 9.974626183044281
$ws.Range("N7").Value = 23.02560710550187

$ws.Range("B8").Value = 17.84730819688782
$ws.Range("D8").Value = 7.632758419142059
$ws.Range("E8").Value = 12.39947351583333
$ws.Range("F8").Value = 39.40921384223044
$ws.Range("G8").Value = 3.707892832446654
$ws.Range("K8").Value = 14.42923233205843
$ws.Range("L8").Value = 10.01541476925895
$ws.Range("N8").Value = 22.97752641257214

$ws.Range("B9").Value = 18.10629746059912
$ws.Range("D9").Value = 7.594454385472201
$ws.Range("E9").Value = 12.43339360309009
$ws.Range("F9").Value = 40.07442143444892
$ws.Range("G9").Value = 3.699483987281567
$ws.Range("K9").Value = 15.10886589470707
$ws.Range("L9").Value = 10.11320795327748
$ws.Range("N9").Value = 22.89796320987545

$ws.Range("B10").Value = 18.31445758747083
$ws.Range("D10").Value = 7.568589060370476
$ws.Range("E10").Value = 12.45799302530456
$ws.Range("F10").Value = 40.60229328232914
$ws.Range("G10").Value = 3.693854743640885
$ws.Range("K10").Value = 15.61020803825244
$ws.Range("L10").Value = 10.19553123183074
$ws.Range("N10").Value = 22.84850362879872

$ws.Range("B11").Value = 18.41274473947837
$ws.Range("D11").Value = 7.557313209170928
$ws.Range("E11").Value = 12.46912265385495
$ws.Range("F11").Value = 40.85024387650425
$ws.Range("G11").Value = 3.691411525805421
$ws.Range("K11").Value = 15.83756559714158
$ws.Range("L11").Value = 10.23516278099208
$ws.Range("N11").Value = 22.82796101065186

$ws.Range("B12").Value = 18.45045404839597
$ws.Range("D12").Value = 7.553113626426165
$ws.Range("E12").Value = 12.47332907416256
$ws.Range("F12").Value = 40.94519842651803
$ws.Range("G12").Value = 3.690503133105341
$ws.Range("K12").Value = 15.92346730482537
$ws.Range("L12").Value = 10.25047519786559
$ws.Range("N12").Value = 22.82046385601656

$ws.Range("B13").Value = 18.44231131168909
$ws.Range("D13").Value = 7.554014956280701
$ws.Range("E13").Value = 12.47242349804388
$ws.Range("F13").Value = 40.92470210571521
$ws.Range("G13").Value = 3.690698026141113
$ws.Range("K13").Value = 15.90497686992632
$ws.Range("L13").Value = 10.24716398724063
$ws.Range("N13").Value = 22.82206595717875

$ws.Range("B14").Value = 18.4158374409615
$ws.Range("D14").Value = 7.556966298334114
$ws.Range("E14").Value = 12.46946887797373
$ws.Range("F14").Value = 40.8580349498768
$ws.Range("G14").Value = 3.691336455666322
$ws.Range("K14").Value = 15.84463717628753
$ws.Range("L14").Value = 10.23641648730435
$ws.Range("N14").Value = 22.82733856283742

$ws.Range("B15").Value = 18.39968444749188
$ws.Range("D15").Value = 7.558783235227345
$ws.Range("E15").Value = 12.46765804751563
$ws.Range("F15").Value = 40.8173356550253
$ws.Range("G15").Value = 3.691729697099556
$ws.Range("K15").Value = 15.80764948842154
$ws.Range("L15").Value = 10.22987275245126
$ws.Range("N15").Value = 22.83060491142048

$ws.Range("B16").Value = 18.30810434911248
$ws.Range("D16").Value = 7.569335816305482
$ws.Range("E16").Value = 12.45726450708672
$ws.Range("F16").Value = 40.5862411034586
$ws.Range("G16").Value = 3.694016770316741
$ws.Range("K16").Value = 15.5953275093951
$ws.Range("L16").Value = 10.19298437262707
$ws.Range("N16").Value = 22.84988555582785

$ws.Range("B17").Value = 18.25282473280473
$ws.Range("D17").Value = 7.575934957119845
$ws.Range("E17").Value = 12.45087327429848
$ws.Range("F17").Value = 40.44643037620701
$ws.Range("G17").Value = 3.69544985033693
$ws.Range("K17").Value = 15.46482719400076
$ws.Range("L17").Value = 10.17090724585239
$ws.Range("N17").Value = 22.86221514468638

$ws.Range("B18").Value = 18.22136931263225
$ws.Range("D18").Value = 7.5797767791186
$ws.Range("E18").Value = 12.44719144927309
$ws.Range("F18").Value = 40.36675485145153
$ws.Range("G18").Value = 3.696285190876949
$ws.Range("K18").Value = 15.38970571817053
$ws.Range("L18").Value = 10.15841504622931
$ws.Range("N18").Value = 22.86949095013068

$ws.Range("B19").Value = 18.21077823421737
$ws.Range("D19").Value = 7.581085489601182
$ws.Range("E19").Value = 12.44594383932283
$ws.Range("F19").Value = 40.33990706088976
$ws.Range("G19").Value = 3.696569927403333
$ws.Range("K19").Value = 15.36426325531763
$ws.Range("L19").Value = 10.15422104480903
$ws.Range("N19").Value = 22.87198602584685

$ws.Range("B20").Value = 18.25867434926592
$ws.Range("D20").Value = 7.575227690215081
$ws.Range("E20").Value = 12.45155422326616
$ws.Range("F20").Value = 40.46123734133819
$ws.Range("G20").Value = 3.695296151487294
$ws.Range("K20").Value = 15.47872618555625
$ws.Range("L20").Value = 10.1732361391164
$ws.Range("N20").Value = 22.86088357637099

$ws.Range("B21").Value = 18.42360037479919
$ws.Range("D21").Value = 7.55609751012287
$ws.Range("E21").Value = 12.47033693722635
$ws.Range("F21").Value = 40.87758842836877
$ws.Range("G21").Value = 3.69114847828786
$ws.Range("K21").Value = 15.86236637615688
$ws.Range("L21").Value = 10.23956509041584
$ws.Range("N21").Value = 22.82578221629326

$ws.Range("B22").Value = 18.53423273974125
$ws.Range("D22").Value = 7.544004731376597
$ws.Range("E22").Value = 12.48256543466291
$ws.Range("F22").Value = 41.15585171265634
$ws.Range("G22").Value = 3.688535614862724
$ws.Range("K22").Value = 16.11192898869757
$ws.Range("L22").Value = 10.28468702495528
$ws.Range("N22").Value = 22.80448475320236

$ws.Range("B23").Value = 18.47493501900966
$ws.Range("D23").Value = 7.550421427990929
$ws.Range("E23").Value = 12.47604296130094
$ws.Range("F23").Value = 41.00679573703823
$ws.Range("G23").Value = 3.689921227111306
$ws.Range("K23").Value = 15.9788689509991
$ws.Range("L23").Value = 10.26044555432422
$ws.Range("N23").Value = 22.81570107870851

$ws.Range("B24").Value = 18.25602872153826
$ws.Range("D24").Value = 7.575547296385612
$ws.Range("E24").Value = 12.45124638944537
$ws.Range("F24").Value = 40.45454091613291
$ws.Range("G24").Value = 3.695365603082632
$ws.Range("K24").Value = 15.47244274063303
$ws.Range("L24").Value = 10.17218262209252
$ws.Range("N24").Value = 22.86148499479005

$ws.Range("B25").Value = 18.03299812267266
$ws.Range("D25").Value = 7.604415686296109
$ws.Range("E25").Value = 12.4242766458668
$ws.Range("F25").Value = 39.887365086027
$ws.Range("G25").Value = 3.701661939582017
$ws.Range("K25").Value = 14.92425040964198
$ws.Range("L25").Value = 10.08488278067921
$ws.Range("N25").Value = 22.91790894001304

Write-Output "Updated loading_percent values for 380 kV case (rows 2-25)"
